$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "Lee's Property"
$ws.Range("G9").Select()
